# Update the cryptos list with fresh price / volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16 and 17 swap their Coin/Link/Price content: WrappedEther moves up
# to rank 16 (pushing ShibaInu to rank 17), each also getting its own new
# Volume(1h) value.
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.652.53"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = "  +1.30%  "

# row => Price (new value, text-forced with a leading apostrophe where the
# text would otherwise be re-interpreted by Excel as a number) and Volume(1h).
$updates = @(
    @{ Row = 2; Price = "59.555.59"; Volume = "  +0.58%  " }
    @{ Row = 3; Price = "2.639.18"; Volume = "  +1.49%  " }
    @{ Row = 4; Price = $null; Volume = "  +0.00%  " }
    @{ Row = 5; Price = "'536.43"; Volume = "  -0.89%  " }
    @{ Row = 6; Price = "'145.53"; Volume = "  +3.17%  " }
    @{ Row = 7; Price = $null; Volume = "  -0.07%  " }
    @{ Row = 8; Price = $null; Volume = "  +1.16%  " }
    @{ Row = 9; Price = "'6.93"; Volume = "  +7.65%  " }
    @{ Row = 10; Price = "'0.101"; Volume = "  -0.99%  " }
    @{ Row = 11; Price = $null; Volume = "  +0.93%  " }
    @{ Row = 12; Price = $null; Volume = "  +0.29%  " }
    @{ Row = 13; Price = "3.104.81"; Volume = "  +1.39%  " }
    @{ Row = 14; Price = "59.448.45"; Volume = "  +0.53%  " }
    @{ Row = 15; Price = "'21.32"; Volume = "  +3.93%  " }
    @{ Row = 18; Price = "'4.48"; Volume = "  +2.73%  " }
    @{ Row = 19; Price = "'338.11"; Volume = "  -1.16%  " }
    @{ Row = 20; Price = "'10.30"; Volume = "  +1.95%  " }
    @{ Row = 21; Price = $null; Volume = "  -2.60%  " }
    @{ Row = 22; Price = $null; Volume = "  -0.02%  " }
    @{ Row = 23; Price = "'66.26"; Volume = "  -1.99%  " }
    @{ Row = 24; Price = $null; Volume = "  +2.35%  " }
    @{ Row = 25; Price = $null; Volume = "  -0.52%  " }
    @{ Row = 26; Price = "'0.999"; Volume = "  -0.07%  " }
    @{ Row = 27; Price = "'7.28"; Volume = "  +1.28%  " }
    @{ Row = 28; Price = "0.0₃0747"; Volume = "  +1.50%  " }
    @{ Row = 30; Price = $null; Volume = "  -2.76%  " }
    @{ Row = 31; Price = "'5.92"; Volume = "  +1.99%  " }
    @{ Row = 32; Price = "'18.83"; Volume = "  +0.69%  " }
    @{ Row = 33; Price = "'150.87"; Volume = "  +0.99%  " }
    @{ Row = 34; Price = $null; Volume = "  +0.67%  " }
    @{ Row = 35; Price = $null; Volume = "  +2.16%  " }
    @{ Row = 36; Price = "'0.839"; Volume = "  +2.96%  " }
    @{ Row = 37; Price = "'0.834"; Volume = "  +0.11%  " }
    @{ Row = 38; Price = $null; Volume = "  -0.91%  " }
    @{ Row = 39; Price = "'3.60"; Volume = "  +1.36%  " }
    @{ Row = 40; Price = "'285.97"; Volume = "  +4.41%  " }
    @{ Row = 41; Price = "'0.999"; Volume = "  -0.04%  " }
    @{ Row = 42; Price = "'0.601"; Volume = "  +0.92%  " }
    @{ Row = 43; Price = $null; Volume = "  +0.04%  " }
    @{ Row = 44; Price = "'0.0537"; Volume = "  +2.49%  " }
    @{ Row = 45; Price = "'19.17"; Volume = "  +2.93%  " }
    @{ Row = 46; Price = "'0.0942"; Volume = "  -1.41%  " }
    @{ Row = 47; Price = $null; Volume = "  +1.48%  " }
    @{ Row = 48; Price = "1.959.42"; Volume = "  +1.04%  " }
    @{ Row = 49; Price = $null; Volume = "  +1.24%  " }
    @{ Row = 50; Price = "'18.37"; Volume = "  -0.32%  " }
    @{ Row = 51; Price = "'111.22"; Volume = "  -0.06%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.Price) {
        $ws.Range("D$row").Value = $u.Price
    }
    $ws.Range("E$row").Value = $u.Volume
}
